$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 307, shifting the existing rows 307-336 down to 310-339.
$ws.Rows.Item(307).Insert()
$ws.Rows.Item(307).Insert()
$ws.Rows.Item(307).Insert()

# Populate the 3 new rows (307-309) with this week's new entries for
# Comercializadora del Agro de Limarí - Alcachofa.

# Row 307
$ws.Cells.Item(307,1).Value = 2
$ws.Cells.Item(307,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(307,3).Value = "Coquimbo"
$ws.Cells.Item(307,4).Value = "5/31/2023"
$ws.Cells.Item(307,5).Value = 4
$ws.Cells.Item(307,6).Value = 100112013
$ws.Cells.Item(307,7).Value = "Alcachofa"
$ws.Cells.Item(307,8).Value = "Argentina(o)"
$ws.Cells.Item(307,9).Value = "Primera"
$ws.Cells.Item(307,10).Value = 700
$ws.Cells.Item(307,11).Value = 9000
$ws.Cells.Item(307,12).Value = 10000
$ws.Cells.Item(307,13).Value = 9500
$ws.Cells.Item(307,14).Value = "$/caja 50 unidades"
$ws.Cells.Item(307,15).Value = "Provincia de Limarí"
$ws.Cells.Item(307,16).Value = 190
$ws.Cells.Item(307,17).Value = 50
$ws.Cells.Item(307,18).Value = "Hortaliza"

# Row 308
$ws.Cells.Item(308,1).Value = 2
$ws.Cells.Item(308,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(308,3).Value = "Coquimbo"
$ws.Cells.Item(308,4).Value = "5/31/2023"
$ws.Cells.Item(308,5).Value = 4
$ws.Cells.Item(308,6).Value = 100112013
$ws.Cells.Item(308,7).Value = "Alcachofa"
$ws.Cells.Item(308,8).Value = "Española"
$ws.Cells.Item(308,9).Value = "Primera"
$ws.Cells.Item(308,10).Value = 800
$ws.Cells.Item(308,11).Value = 13000
$ws.Cells.Item(308,12).Value = 14000
$ws.Cells.Item(308,13).Value = 13500
$ws.Cells.Item(308,14).Value = "$/caja 30 unidades"
$ws.Cells.Item(308,15).Value = "Provincia de Limarí"
$ws.Cells.Item(308,16).Value = 450
$ws.Cells.Item(308,17).Value = 30
$ws.Cells.Item(308,18).Value = "Hortaliza"

# Row 309
$ws.Cells.Item(309,1).Value = 2
$ws.Cells.Item(309,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(309,3).Value = "Coquimbo"
$ws.Cells.Item(309,4).Value = "5/31/2023"
$ws.Cells.Item(309,5).Value = 4
$ws.Cells.Item(309,6).Value = 100112013
$ws.Cells.Item(309,7).Value = "Alcachofa"
$ws.Cells.Item(309,8).Value = "Madrigal"
$ws.Cells.Item(309,9).Value = "Primera"
$ws.Cells.Item(309,10).Value = 500
$ws.Cells.Item(309,11).Value = 13000
$ws.Cells.Item(309,12).Value = 14000
$ws.Cells.Item(309,13).Value = 13500
$ws.Cells.Item(309,14).Value = "$/caja 40 unidades"
$ws.Cells.Item(309,15).Value = "Provincia de Limarí"
$ws.Cells.Item(309,16).Value = 338
$ws.Cells.Item(309,17).Value = 40
$ws.Cells.Item(309,18).Value = "Hortaliza"
